$d = $word.ActiveDocument

# --- 1. Remove the two standalone paragraphs "Change structure," and
#        "More information about static" (whole paragraphs, incl. marks) ---
$pChange = $null
$pMore = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -eq "Change structure,`r") {
        $pChange = $p
    }
    elseif ($t -eq "More information about static`r") {
        $pMore = $p
    }
}
if (($pChange -ne $null) -and ($pMore -ne $null)) {
    $delRange = $d.Range($pChange.Range.Start, $pMore.Range.End)
    $delRange.Delete()
}

# --- 2. Collapse the four runs that make up the "We want to ..." paragraph
#        into a single run containing the same combined text. ---
$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$fullText = "We want to " + $quoteOpen + "Use wireframe" + $quoteClose + " to make an attractive way for the user to educate himself. Only a little information for every kind of competition, explaining of the tools themselves."

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($fullText + "`r")) {
        $target = $p
    }
}
if ($target -ne $null) {
    $searchRange = $target.Range
    $f = $searchRange.Find
    $f.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, $fullText, 2) | Out-Null
}
